$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7141596666666666
$ws.Range("H2").Value = 2.142479
$ws.Range("O2").Value = 0.9886227745742286
$ws.Range("P2").Value = 0.9886227745742288
$ws.Range("Q2").Value = 0.4069893577447777
$ws.Range("R2").Value = 3.662904219703
$ws.Range("S2").Value = 0.9886227745742286
$ws.Range("T2").Value = 0.9886227745742288

# Add new row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Insl3"
$ws.Range("C3").Value = "Rxfp2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7141596666666666
$ws.Range("H3").Value = 2.142479
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.006558333333333333
$ws.Range("N3").Value = 0.019675
$ws.Range("O3").Value = 0.01137722542577134
$ws.Range("P3").Value = 0.01137722542577134
$ws.Range("Q3").Value = 0.004683697147222222
$ws.Range("R3").Value = 0.04215327432499999
$ws.Range("S3").Value = 0.01137722542577134
$ws.Range("T3").Value = 0.01137722542577134
